$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.418.15'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.619.06'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.72%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '203.38'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +10.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '568.25'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.615.54'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.71%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.55%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.677'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '61.30'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +16.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.152'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +4.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000288'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +10.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.08'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.198.20'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.622.58'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.62%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.92%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '68.297.54'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.41'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '404.32'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.15'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +17.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.18'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.57'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.65'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.77%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +10.42%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.15'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +15.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.40'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +5.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.67'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '675.44'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +9.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.31'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.34%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '63.73'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.31'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.69%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.81%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0778'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.22%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.21'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +14.21%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.261.16'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +9.61%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.76'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +10.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.01'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +29.39%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.82'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +12.90%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0420'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.91'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.21%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.09'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.51%  '
